$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the latitude value in B2 (simulation input changed)
$ws.Range("B2").Value = 79.81

# Update selection to reflect where the user was working
$ws.Range("H19").Select()
